$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "want to go" counts (column F)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 1775
$wsExhibit.Range("F6").Value = 270

# Sheet "全部类型" (all types) - same events appear here, update matching cells
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1775
$wsAll.Range("F7").Value = 270
